# Nueva tarea aleatorizar botones talk to
#
# Locate the empty "Prrafodelista" paragraph that currently only holds the
# _GoBack bookmark (right after "Repetir esto nueve veces (Miguel y Juan
# Camilo)"). Turn it into a new list item (continuing the same numbered
# list as its neighbours, numId=6) with the new task text, then leave a
# fresh empty "Prrafodelista" paragraph (no numbering) right after it --
# mirroring the paragraph that used to hold the bookmark.

$d = $word.ActiveDocument

$anchorText = "Repetir esto nueve veces (Miguel y Juan Camilo)"
$prevIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $candText = $cand.Range.Text
    if ($candText.TrimEnd("`r") -eq $anchorText) {
        $prevIndex = $i
        break
    }
}

$targetIndex = $prevIndex + 1
$prev = $d.Paragraphs.Item($prevIndex)
$target = $d.Paragraphs.Item($targetIndex)

# Continue the preceding numbered list (numId=6) on the target paragraph.
$listTemplate = $prev.Range.ListFormat.ListTemplate
$target.Range.ListFormat.ApplyListTemplate($listTemplate, $true)

# Insert the task text (with embedded proofErr spell-check markers) right
# before the bookmark, i.e. at the very start of the (still empty) paragraph.
$target = $d.Paragraphs.Item($targetIndex)
$insertAt = $target.Range.Start
$insertRange = $d.Range($insertAt, $insertAt)
$taskXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:r><w:t>Hacer que los botones dentro de la ventana de &#8220;</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Talk</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> to </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>partner</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t>&#8221; aparezcan de manera aleatoria</w:t></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($taskXml)

# Re-fetch the (now non-empty) paragraph and append a new blank
# "Prrafodelista" paragraph right after it, without any numbering.
$target = $d.Paragraphs.Item($targetIndex)
$afterAt = $target.Range.End
$afterRange = $d.Range($afterAt, $afterAt)
$blankXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/></w:pPr></w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$afterRange.InsertXML($blankXml)
